$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: reorder "Recorded By" list - move "System" earlier in the sequence
$ws.Range("G2").Value = "Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

# G9: swap order of the two recorders
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# L10 & S15: percentage text cells ("10.2%" -> "18.2%").
# These are stored as plain text (General-formatted cells), so a direct
# Value assignment would make Excel auto-convert the string into a real
# percentage number (and change the cell's number format / style).
# Route the new text through a formula + Copy/PasteSpecial(values) so the
# result lands back as literal text in the original style.
$ws.Range("L10").Formula = "=""18.2%"""
$ws.Range("L10").Copy()
$ws.Range("L10").PasteSpecial(-4163)

$ws.Range("S15").Formula = "=""18.2%"""
$ws.Range("S15").Copy()
$ws.Range("S15").PasteSpecial(-4163)

# G28: append an additional recorder email
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

# H28: attendance count 6/251 -> 66/251
$ws.Range("H28").Value = "66/251"

$excel.CutCopyMode = 0
